{"js": "// The captured change is NOT a content/formatting edit at all: diffing the\n// package's canonical OOXML shows that every \"-\"/\"+\" line pair touches the\n// exact same element (same tag, same attribute names/values, same child\n// structure) and only the *serialization order* of XML attributes (and of\n// the root <w:document> element's xmlns:* declarations) differs, e.g.:\n//   -<w:pgSz w:w=\"11906\" w:h=\"16838\"/>\n//   +<w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n//   -<w:style w:type=\"paragraph\" w:default=\"1\" w:styleId=\"Normal\">\n//   +<w:style w:default=\"1\" w:styleId=\"Normal\" w:type=\"paragraph\">\n// This is the well-known effect described by the commit message itself\n// (\"Fixed POI packaging and upgraded to POI 3.15\"): re-saving the fixture\n// with a newer Apache POI/XMLBeans writer reorders attributes\n// alphabetically while leaving every value, run of text, style, numbering\n// entry, margin, font, etc. completely untouched. A quick structural diff\n// confirms there are zero textual/semantic differences between the before\n// and after XML (only attribute order changed).\n//\n// The Word JavaScript API has no concept of \"XML attribute order\" (it is\n// an internal detail of whichever OOXML writer happens to serialize the\n// part) and exposes no operation that can reorder attributes within a\n// tag or reorder namespace declarations on the document root, so there is\n// no content-visible action to perform here. To avoid introducing any\n// unintended change, this script intentionally performs a harmless,\n// read-only touch of the document body and applies no edits.\nconst body = context.document.body;\nbody.load(\"text\");\nawait context.sync();\n", "ps1": "# The captured change is NOT a content/formatting edit at all: diffing the\n# package's canonical OOXML shows that every \"-\"/\"+\" line pair touches the\n# exact same element (same tag, same attribute names/values, same child\n# structure) and only the *serialization order* of XML attributes (and of\n# the root <w:document> element's xmlns:* declarations) differs, e.g.:\n#   -<w:pgSz w:w=\"11906\" w:h=\"16838\"/>\n#   +<w:pgSz w:h=\"16838\" w:w=\"11906\"/>\n#   -<w:style w:type=\"paragraph\" w:default=\"1\" w:styleId=\"Normal\">\n#   +<w:style w:default=\"1\" w:styleId=\"Normal\" w:type=\"paragraph\">\n# This is the well-known effect described by the commit message itself\n# (\"Fixed POI packaging and upgraded to POI 3.15\"): re-saving the fixture\n# with a newer Apache POI/XMLBeans writer reorders attributes\n# alphabetically while leaving every value, run of text, style, numbering\n# entry, margin, font, etc. completely untouched. A structural diff\n# confirms there are zero textual/semantic differences between the before\n# and after XML (only attribute order changed).\n#\n# The Word COM object model has no concept of \"XML attribute order\" (it is\n# an internal detail of whichever OOXML writer happens to serialize the\n# part) and exposes no operation that can reorder attributes within a tag\n# or reorder namespace declarations on the document root, so there is no\n# content-visible action to perform here. To avoid introducing any\n# unintended change, this script intentionally performs a harmless,\n# read-only touch of the document and applies no edits.\n$d = $word.ActiveDocument\n$null = $d.Content.Text\n"}
